# Update dSF (column F) values to reflect a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = 0
    4  = -5
    8  = -1
    9  = 3
    10 = 0
    14 = 2
    15 = 1
    17 = 2
    20 = 0
    23 = 3
    25 = 4
    34 = -4
    35 = -2
    37 = -2
    40 = 3
    45 = 1
    46 = -4
    48 = 1
    52 = 5
    53 = 3
    62 = -4
    74 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
